# "create entity in progress" - refresh the CreateEntity_OOFS_NewValues test
# fixture with values from a newer automation run (2024-04-06 run instead of
# the 2024-03-23..28 runs baked into the previous snapshot).

$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Add New") ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G2").Value  = "5803514002"
$ws1.Range("O2").Value  = "2024-04-06"
$ws1.Range("Q2").Value  = "2024-04-06 02:50:13 PM"
$ws1.Range("AD2").Value = "2024-04-06"
$ws1.Range("AF2").Value = "2402685023"
$ws1.Range("AV2").Value = "2392196750"
$ws1.Range("AZ2").Value = "7013689143"
$ws1.Range("BB2").Value = "CT: Sat, Apr 06, 2024 at 2:52 PM"

# --- Sheet2 ("Summary Add") ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("G2").Value  = "5803514002"
$ws2.Range("O2").Value  = "2024-04-06"
$ws2.Range("Q2").Value  = "2024-04-06 02:56:00 PM"
$ws2.Range("AD2").Value = "2024-04-06"
$ws2.Range("AF2").Value = "2402685023"
$ws2.Range("AV2").Value = "2392196750"
$ws2.Range("AZ2").Value = "7013689143"
$ws2.Range("BB2").Value = "CT: Sat, Apr 06, 2024 at 3:01 PM"

# --- Sheet3 ("Duplicate") ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("G2").Value  = "5803514002"
$ws3.Range("AF2").Value = "2402685023"
$ws3.Range("AV2").Value = "2392196750"
$ws3.Range("AZ2").Value = "7013689143"

# --- Sheet4 ("Edit Record") ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("G2").Value  = "5803514002"
$ws4.Range("AF2").Value = "2402685023"
$ws4.Range("AV2").Value = "2392196750"
$ws4.Range("AZ2").Value = "7013689143"

# --- Sheet5 ("Default Success") ---
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("AI2").Value = "Create Entity - CreateEntityTarget"
$ws5.Range("AJ2").Value = "tskOOFS_CEOnly"
